$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct two previously mis-marked attendance cells (column AF, date 2025-08-21)
# Row 13 (Rayane Chayebi): was "B" (Blessure), should be "RH"
# Row 18 (Emmanuel Valey): was "P" (Présent), should be "RH"
$ws.Range("AF13").Value = "RH"
$ws.Range("AF18").Value = "RH"

# Add the new training-day column AG for 2025-08-23 (serial date 45892)
$ws.Range("AG1").Value = 45892

# Attendance codes for the new date, one per player row
$values = @{
  2  = "P"
  3  = "P"
  4  = "P"
  5  = "P"
  6  = "P"
  7  = "P"
  8  = "R"
  9  = "P"
  10 = "P"
  11 = "P"
  12 = "P"
  13 = "P"
  14 = "P"
  15 = "P"
  16 = "P"
  17 = "B"
  18 = "RH"
  19 = "RH"
  20 = "P"
  21 = "R"
  22 = "P"
  23 = "R"
  24 = "P"
  25 = "P"
  26 = "P"
  27 = "RH"
}

foreach ($r in $values.Keys) {
  $ws.Range("AG" + $r).Value = $values[$r]
}

# Copy the formatting (style, number format) from column AF onto the new column AG
$src = $ws.Range("AF1:AF27")
$src.Copy()
$ws.Range("AG1:AG27").PasteSpecial(-4122)

# Move the active cell selection one column over, following the newly added column
$ws.Range("AI24").Select() | Out-Null
